# Apply the cryptos list update (Wed Nov  1 06:45:14 UTC 2023 run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.446.93"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.804.12"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'224.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "'0.591"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'38.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.19%  "
$ws.Range("D9").Value = "'0.287"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.60%  "
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").Value = "'0.0974"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "2.064.32"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "'11.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.95%  "
$ws.Range("D14").Value = "1.804.28"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "34.421.00"
$ws.Range("D16").Value = "'0.626"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "'67.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").Value = "'241.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").Value = "0.0₃0768"
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("D21").Value = "'11.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'4.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.98%  "
$ws.Range("D24").Value = "'2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("D25").Value = "'170.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("D27").Value = "'17.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "'1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'3.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").Value = "'3.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("D33").Value = "'0.0513"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").Value = "1.324.20"
$ws.Range("E35").Value = "  -5.35%  "
$ws.Range("E36").Value = "  -5.01%  "
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").Value = "'82.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'2.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.56%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("D44").Value = "'0.941"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").Value = "'13.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").Value = "1.964.81"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  -4.98%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'101.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").Value = "0.0₆0120"
$ws.Range("E51").Value = "  -6.32%  "
